# Update mods data [2026-02-15 15:14:32]
# Append a new data row (row 97) to the ModCounts sheet:
#   2026/02/15 | 逃离鸭科夫 | 1208

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 97

# Column A holds a date-like string ("2026/02/15") that must stay a literal
# text value (matching the sheet's existing inline-string date cells) rather
# than being auto-converted into a date serial number by Excel's input
# parsing. Briefly mark the cell as Text before assigning the value, then
# strip the resulting explicit format so the cell falls back to the same
# look (General number format, centered alignment) as the rest of the table.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2026/02/15"
$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1208

$rowRange = $ws.Range("A$newRow`:C$newRow")
$rowRange.ClearFormats()
$rowRange.HorizontalAlignment = -4108
$rowRange.VerticalAlignment = -4108
